$d = $word.ActiveDocument

# Locate the paragraph that ends with "const <- matrix(..., byrow = TRUE)"
$constPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*byrow = TRUE)*") {
        $constPara = $p
    }
}

$r = $constPara.Range

# --- Append two line breaks followed by "const" (styled NormalTok) to the
#     end of the existing paragraph (before its paragraph mark). ---
$ip = $d.Range($r.End - 1, $r.End - 1)
$ip.InsertBreak(6)

$ip = $d.Range($r.End - 1, $r.End - 1)
$ip.InsertBreak(6)

$ip = $d.Range($r.End - 1, $r.End - 1)
$constStart = $ip.Start
$ip.InsertAfter("const")
$constRange = $d.Range($constStart, $ip.End)
$constRange.Style = "NormalTok"

# --- Insert a brand-new SourceCode paragraph right after it, holding the
#     printed constraint matrix (each line styled VerbatimChar, separated
#     by manual line breaks). ---
$ip = $d.Range($r.End - 1, $r.End - 1)
$ip.InsertParagraphAfter()

$newPara = $constPara.Next()
$newRange = $newPara.Range

$lines = @(
    "##      [,1] [,2] [,3] [,4] [,5] [,6] [,7]",
    "## [1,]    0    1    1    1    1    1    0",
    "## [2,]    0    0    1    1    1    1    1",
    "## [3,]    1    0    0    1    1    1    1",
    "## [4,]    1    1    0    0    1    1    1",
    "## [5,]    1    1    1    0    0    1    1",
    "## [6,]    1    1    1    1    0    0    1",
    "## [7,]    1    1    1    1    1    0    0"
)

for ($i = 0; $i -lt $lines.Length; $i++) {
    $ip = $d.Range($newRange.End - 1, $newRange.End - 1)
    $lineStart = $ip.Start
    $ip.InsertAfter($lines[$i])
    $lineRange = $d.Range($lineStart, $ip.End)
    $lineRange.Style = "VerbatimChar"

    if ($i -lt $lines.Length - 1) {
        $ip = $d.Range($newRange.End - 1, $newRange.End - 1)
        $ip.InsertBreak(6)
    }
}

Write-Output "Edit complete"
